$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.316.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -15.54%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.274.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -21.66%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '432.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -17.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -19.35%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.452'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -16.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.270.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -22.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.09'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -16.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -21.46%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.295'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -17.64%  '

$ws.Range('E13').Value = '  -6.85%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.653.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -22.21%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '51.445.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -15.27%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -19.27%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000112'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -19.95%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.272.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -21.96%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -21.19%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '289.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -17.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.45%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -27.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -24.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.995'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '52.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -19.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.359'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -20.48%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.326.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -23.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.138'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -21.59%  '

$ws.Range('E30').Value = '  -0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -16.10%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '141.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₃0617'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -27.69%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '16.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -16.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -22.75%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -18.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.994'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.27'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -25.45%  '

$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.950'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -20.36%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '31.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -15.91%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.732'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -26.17%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.86%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.548'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -15.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -18.48%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0483'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -16.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.845.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -19.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -24.95%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0197'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -16.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0781'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.69%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '15.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -25.14%  '

$ws.Range('B51').Value = 'ZEEBU'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.66%  '
